# Incorporacion de Logica de Anexos Postgres
#
# Row 35 ("Alvear" / sigehoslgc_salvear / id 1644) in the "Tabla1" table was
# an exact duplicate of row 3. Remove that duplicate row; everything below
# shifts up by one (table range, dimension and autofilter follow suit
# automatically), and the view/selection is reset to the top of the sheet
# with B5 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(35).Delete()

$ws.Range("B5").Select()
